$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.864.34"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "3.759.92"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "645.83"
$ws.Range("E5").Value = "  +2.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.18"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").Value = "3.758.73"
$ws.Range("E7").Value = "  -1.27%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.455"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.90"
$ws.Range("E12").Value = "  +4.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  -4.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.73"
$ws.Range("E14").Value = "  -3.21%  "

$ws.Range("D15").Value = "4.399.15"
$ws.Range("E15").Value = "  -1.11%  "

$ws.Range("D16").Value = "3.761.06"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "68.894.96"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.60"
$ws.Range("E18").Value = "  -1.89%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.99"
$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.96"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.52"
$ws.Range("E22").Value = "  -1.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.702"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000143"
$ws.Range("E24").Value = "  -4.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.74"
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.15"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.09"
$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "3.913.65"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").Value = "  +1.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.09"
$ws.Range("E33").Value = "  -2.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.41"
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.173"
$ws.Range("E35").Value = "  +15.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "3.720.10"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.77"
$ws.Range("E38").Value = "  -2.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("E39").Value = "  -2.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.75"
$ws.Range("E40").Value = "  -2.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.21"
$ws.Range("E41").Value = "  -6.75%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.955"
$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.86"
$ws.Range("E45").Value = "  +4.10%  "

$ws.Range("E46").Value = "  +3.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.04"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.29"
$ws.Range("E48").Value = "  +0.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.40"
$ws.Range("E49").Value = "  -1.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.294"
$ws.Range("E50").Value = "  -1.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.34"
$ws.Range("E51").Value = "  -1.15%  "
